# Commit: "Create contentTable fixes #1"
# とりあえず、テーブルを表示させるところまで編集した
#
# Turns the old "Docomo" table on the first sheet into the new "appleID"
# table: rename the sheet tab, relabel its header content, and make sure
# the sheet shows as a plain (non-frozen) table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the first sheet from "シート1" to "appleID"
$ws.Name = "appleID"

# The cell that used to read "Docomo" now reads "appleID"
$ws.Range("B2").Value = "appleID"

# Make sure the table is displayed without a frozen header pane/split
$ws.Activate()
if ($excel.ActiveWindow.FreezePanes) {
    $excel.ActiveWindow.FreezePanes = $false
}
$excel.ActiveWindow.Split = $false
